$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 595.2174
$ws.Range("J17").Value = 601.0448
$ws.Range("L17").Value = 1803.1344
$ws.Range("N17").Value = -2139.1344

# Row 116
$ws.Range("H116").Value = 15382110
$ws.Range("I116").Value = 27675198
$ws.Range("K116").Value = 27675198
$ws.Range("M116").Value = -27671756

# Row 132
$ws.Range("H132").Value = 595027.75
$ws.Range("I132").Value = 715916.3
$ws.Range("J132").Value = 81251.5
$ws.Range("K132").Value = 2147748.9
$ws.Range("L132").Value = 243754.5
$ws.Range("M132").Value = -2145218.9
$ws.Range("N132").Value = -248814.5

# Row 137
$ws.Range("H137").Value = 125004150
$ws.Range("I137").Value = 250003000
$ws.Range("K137").Value = 750009000
$ws.Range("M137").Value = -750006450

# Row 138
$ws.Range("H138").Value = 6051209.5
$ws.Range("I138").Value = 1539011.4
$ws.Range("J138").Value = 7465480.5
$ws.Range("K138").Value = 4617034.199999999
$ws.Range("L138").Value = 22396441.5
$ws.Range("M138").Value = -4611894.199999999
$ws.Range("N138").Value = -22406721.5

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1199.875
$ws.Range("I45").Value = 1119.8
$ws.Range("K45").Value = 1119.8
$ws.Range("M45").Value = -742.8

# Row 132
$ws.Range("H132").Value = 3223.4285
$ws.Range("I132").Value = 2858.75
$ws.Range("K132").Value = 8576.25
$ws.Range("M132").Value = -6046.25

# Row 133
$ws.Range("H133").Value = 30250
$ws.Range("J133").Value = 30250
$ws.Range("L133").Value = 30250
$ws.Range("N133").Value = -35310

$ws = $wb.Worksheets.Item("BSM")
# Row 19
$ws.Range("H19").Value = 50000
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

# Row 105
$ws.Range("H105").Value = 2553.6206
$ws.Range("I105").Value = 2410.2
$ws.Range("J105").Value = 3450
$ws.Range("K105").Value = 2410.2
$ws.Range("L105").Value = 3450
$ws.Range("M105").Value = -663.1999999999998
$ws.Range("N105").Value = -6944

# Row 107
$ws.Range("H107").Value = 527.3570999999999
$ws.Range("I107").Value = 453
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 453
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 1467
$ws.Range("N107").Value = -4640

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1013
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1013
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 1013
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1587

# Row 31
$ws.Range("H31").Value = 2768.1428
$ws.Range("I31").Value = 1536.9286
$ws.Range("J31").Value = 5230.5713
$ws.Range("K31").Value = 1536.9286
$ws.Range("L31").Value = 5230.5713
$ws.Range("M31").Value = -1241.9286
$ws.Range("N31").Value = -5820.5713

# Row 34
$ws.Range("H34").Value = 2768.1428
$ws.Range("I34").Value = 1536.9286
$ws.Range("J34").Value = 5230.5713
$ws.Range("K34").Value = 1536.9286
$ws.Range("L34").Value = 5230.5713
$ws.Range("M34").Value = -1334.9286
$ws.Range("N34").Value = -5634.5713

# Row 107
$ws.Range("H107").Value = 636.5
$ws.Range("I107").Value = 323
$ws.Range("K107").Value = 323
$ws.Range("M107").Value = 1597

# Row 113
$ws.Range("H113").Value = 1013
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1013
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1013
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -5353

# Row 132
$ws.Range("H132").Value = 2239.121
$ws.Range("I132").Value = 1593.3
$ws.Range("J132").Value = 3232.6924
$ws.Range("K132").Value = 4779.9
$ws.Range("L132").Value = 9698.0772
$ws.Range("M132").Value = -2249.9
$ws.Range("N132").Value = -14758.0772

# Row 134
$ws.Range("H134").Value = 2426.647
$ws.Range("I134").Value = 1426.12
$ws.Range("J134").Value = 5205.8887
$ws.Range("K134").Value = 4278.36
$ws.Range("L134").Value = 15617.6661
$ws.Range("M134").Value = -1743.36
$ws.Range("N134").Value = -20687.6661

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1702.5
$ws.Range("I5").Value = 849
$ws.Range("J5").Value = 2745.6667
$ws.Range("K5").Value = 2547
$ws.Range("L5").Value = 8237.000100000001
$ws.Range("M5").Value = -2435
$ws.Range("N5").Value = -8461.000100000001

# Row 122
$ws.Range("H122").Value = 724.2105
$ws.Range("J122").Value = 962.3333
$ws.Range("L122").Value = 8660.9997
$ws.Range("N122").Value = -13560.9997

# Row 135
$ws.Range("H135").Value = 1702.5
$ws.Range("I135").Value = 849
$ws.Range("J135").Value = 2745.6667
$ws.Range("K135").Value = 7641
$ws.Range("L135").Value = 24711.0003
$ws.Range("M135").Value = -5106
$ws.Range("N135").Value = -29781.0003

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6440.75
$ws.Range("I70").Value = 6806.625
$ws.Range("J70").Value = 4977.25
$ws.Range("K70").Value = 6806.625
$ws.Range("L70").Value = 4977.25
$ws.Range("M70").Value = -6536.625
$ws.Range("N70").Value = -5517.25

# Row 73
$ws.Range("H73").Value = 6440.75
$ws.Range("I73").Value = 6806.625
$ws.Range("J73").Value = 4977.25
$ws.Range("K73").Value = 6806.625
$ws.Range("L73").Value = 4977.25
$ws.Range("M73").Value = -5870.625
$ws.Range("N73").Value = -6849.25

# Row 103
$ws.Range("H103").Value = 18700
$ws.Range("J103").Value = 18700
$ws.Range("L103").Value = 18700
$ws.Range("N103").Value = -21044

# Row 126
$ws.Range("H126").Value = 2886.95
$ws.Range("I126").Value = 2300
$ws.Range("J126").Value = 2990.5293
$ws.Range("K126").Value = 6900
$ws.Range("L126").Value = 8971.5879
$ws.Range("M126").Value = -4430
$ws.Range("N126").Value = -13911.5879

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 3855.2856
$ws.Range("I136").Value = 2787.9
$ws.Range("J136").Value = 4825.636
$ws.Range("K136").Value = 8363.700000000001
$ws.Range("L136").Value = 14476.908
$ws.Range("M136").Value = -5813.700000000001
$ws.Range("N136").Value = -19576.908

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 347.875
$ws.Range("I113").Value = 196
$ws.Range("J113").Value = 398.5
$ws.Range("K113").Value = 588
$ws.Range("L113").Value = 1195.5
$ws.Range("M113").Value = 1582
$ws.Range("N113").Value = -5535.5
